$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = 'MAE [$COP/kWh]'
$ws.Range("K1").Value = 'MSE [$COP/kWh]'
$ws.Range("L1").Value = 'RMSE [$COP/kWh]'
$ws.Range("M1").Value = 'MAPE [%]'

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = 30
$ws.Range("G2").Value = '<keras.src.optimizers.adam.Adam object at 0x000001D1FF045090>'
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 101.5419890534182
$ws.Range("K2").Value = 13366.43791780146
$ws.Range("L2").Value = 115.6133120267794
$ws.Range("M2").Value = 58.0289920242638
